$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of cell updates: (cell reference, new text value) extracted from the
# updated cryptocurrency price/volume snapshot.
$updates = @(
    @{Cell="D2"; Value="302.32"},
    @{Cell="E2"; Value="-5.17%"},
    @{Cell="D3"; Value="35.13"},
    @{Cell="E3"; Value="-2.59%"},
    @{Cell="D4"; Value="5.046"},
    @{Cell="E4"; Value="-1.58%"},
    @{Cell="D5"; Value="0.07974"},
    @{Cell="E5"; Value="-2.66%"},
    @{Cell="D6"; Value="1.911"},
    @{Cell="E6"; Value="-10.52%"},
    @{Cell="D7"; Value="4.038"},
    @{Cell="E7"; Value="-2.33%"},
    @{Cell="D8"; Value="7.732"},
    @{Cell="E8"; Value="-3.46%"},
    @{Cell="D9"; Value="2.953"},
    @{Cell="E9"; Value="5.44%"},
    @{Cell="D10"; Value="0.9230"},
    @{Cell="E10"; Value="-0.32%"},
    @{Cell="D11"; Value="0.1203"},
    @{Cell="E11"; Value="19.58%"},
    @{Cell="D12"; Value="0.1841"},
    @{Cell="E12"; Value="-2.97%"},
    @{Cell="D13"; Value="0.09363"},
    @{Cell="E13"; Value="1.55%"},
    @{Cell="D14"; Value="0.03526"},
    @{Cell="E14"; Value="-2.02%"},
    @{Cell="D15"; Value="0.09852"},
    @{Cell="E15"; Value="-0.69%"},
    @{Cell="D16"; Value="0.001399"},
    @{Cell="E16"; Value="-2.33%"},
    @{Cell="D17"; Value="0.005859"},
    @{Cell="E17"; Value="3.21%"},
    @{Cell="D18"; Value="3.496"},
    @{Cell="E18"; Value="0.93%"},
    @{Cell="D19"; Value="0.3445"},
    @{Cell="E19"; Value="2.12%"},
    @{Cell="E20"; Value="-2.94%"},
    @{Cell="D21"; Value="5.055"},
    @{Cell="E21"; Value="-0.45%"},
    @{Cell="D22"; Value="0.2465"},
    @{Cell="E22"; Value="12.63%"},
    @{Cell="D23"; Value="0.04489"},
    @{Cell="E23"; Value="-2.33%"},
    @{Cell="D24"; Value="0.001214"},
    @{Cell="E24"; Value="-2.31%"},
    @{Cell="D25"; Value="0.004569"},
    @{Cell="E25"; Value="-3.54%"},
    @{Cell="D26"; Value="0.0001250"},
    @{Cell="E26"; Value="-3.87%"},
    @{Cell="E27"; Value="-6.85%"},
    @{Cell="D39"; Value="0.01913"},
    @{Cell="E39"; Value="-4.85%"},
    @{Cell="D40"; Value="0.04735"},
    @{Cell="E40"; Value="-5.11%"},
    @{Cell="D41"; Value="0.007597"},
    @{Cell="E41"; Value="-1.68%"},
    @{Cell="D42"; Value="0.009554"},
    @{Cell="E42"; Value="22.31%"},
    @{Cell="D43"; Value="0.1325"},
    @{Cell="E43"; Value="-5.37%"},
    @{Cell="D44"; Value="0.002110"},
    @{Cell="E44"; Value="-0.96%"},
    @{Cell="D45"; Value="0.01114"},
    @{Cell="E45"; Value="-7.19%"},
    @{Cell="D46"; Value="0.00006273"},
    @{Cell="E46"; Value="-2.66%"},
    @{Cell="D47"; Value="0.00000000750"},
    @{Cell="E47"; Value="-0.03%"},
    @{Cell="E49"; Value="-31.37%"},
    @{Cell="D50"; Value="0.00002100"},
    @{Cell="E50"; Value="-0.03%"},
    @{Cell="D51"; Value="0.0002000"},
    @{Cell="E51"; Value="-0.03%"}
)

foreach ($item in $updates) {
    $cellRef = $item.Cell
    $val = $item.Value
    $range = $ws.Range($cellRef)
    # Build a formula that evaluates to the literal text so the cell keeps
    # its original "text" storage (values like "302.32" or "-5.17%" must not
    # be auto-converted to numeric/percentage types).
    $escapedVal = $val.Replace('"', '""')
    $range.Formula = '="' + $escapedVal + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = 0
